# Update the "想去人数" (column F) counts on the 展览 and 全部类型 sheets.
# These two sheets list (mostly) the same events, so the same new values
# are applied to both, using each sheet's own row numbering.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> hashtable of row -> new value for column F
$updates = @{
    "展览"   = @{ 4 = 88; 6 = 24; 8 = 8045; 9 = 758; 10 = 240; 11 = 1099; 12 = 785; 13 = 34; 15 = 202; 16 = 54; 17 = 51; 18 = 209; 19 = 854 }
    "全部类型" = @{ 4 = 88; 6 = 24; 9 = 8045; 10 = 758; 11 = 240; 12 = 1099; 13 = 785; 14 = 34; 16 = 202; 17 = 54; 18 = 51; 19 = 209; 20 = 854 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsForSheet = $updates[$sheetName]
    foreach ($row in $rowsForSheet.Keys) {
        $ws.Range("F$row").Value = $rowsForSheet[$row]
    }
}
